$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCL questions prod")
Write-Host $ws.Name
